$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The transmitter incidence input ("Tx_th", zenith angle) was converted to
# an elevation input ("Tx_el") - update the column header text accordingly.
$ws.Range("A1").Value = "Tx_el (deg)"

# Update the sheet's stored selection from G2 to B10.
$ws.Range("B10").Select() | Out-Null
